$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 corresponds to "QUICK NAIL LOTION" - its "current balance" cell (H9)
# changes from "-23:0" to "0:0", matching the other rows' balance text
# ("-1:0" -> "0:0") that is shared with H7/H8.
$ws.Range("H7").Value = "0:0"
$ws.Range("H8").Value = "0:0"
$ws.Range("H9").Value = "0:0"
